$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the new value "Regiane" (adds a new shared string entry)
$ws.Range("A2").Value = "Regiane"

# Move the selection from D2 to A2
$ws.Range("A2").Select()
